$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-21, columns C (chromosome string), D (Maximo), E (Minimo), F (Promedio)
$data = @(
    @{ Row = 2;  C = "111101111010101101010110110011"; D = 0.9359753248780518; E = 0.04288710317249519; F = 0.5449269110666202 },
    @{ Row = 3;  C = "111101111010101101010110110011"; D = 0.9359753248780518; E = 0.4895102378827652;  F = 0.6981919740518754 },
    @{ Row = 4;  C = "111101111010101101010110110111"; D = 0.9359753320861769; E = 0.5273611445630106;  F = 0.7745138896389954 },
    @{ Row = 5;  C = "111101111010101101010110110111"; D = 0.9359753320861769; E = 0.5273611445630106;  F = 0.8227086109155509 },
    @{ Row = 6;  C = "111101111010101101010110110111"; D = 0.9359753320861769; E = 0.5273611445630106;  F = 0.8349721765401965 },
    @{ Row = 7;  C = "111110101001010011110001000111"; D = 0.9581182957678596; E = 0.8188921376246626;  F = 0.8975775265316546 },
    @{ Row = 8;  C = "111110101001010011110001000111"; D = 0.9581182957678596; E = 0.8188917297198883;  F = 0.8975801087417248 },
    @{ Row = 9;  C = "111110101001010011110001000111"; D = 0.9581182957678596; E = 0.8188923500048747;  F = 0.9070482143204645 },
    @{ Row = 10; C = "111110101001010011110001000111"; D = 0.9581182957678596; E = 0.8396701490447454;  F = 0.9111820498120731 },
    @{ Row = 11; C = "111110101001010011110110110111"; D = 0.9581189667138894; E = 0.8396701490447454;  F = 0.9134624661246699 },
    @{ Row = 12; C = "111110101001010011110110110111"; D = 0.9581189667138894; E = 0.8396701490447454;  F = 0.9038980650874031 },
    @{ Row = 13; C = "111110101001010011110110110111"; D = 0.9581189667138894; E = 0.8396701490447454;  F = 0.8897727670078919 },
    @{ Row = 14; C = "111110101001010011110110110111"; D = 0.9581189667138894; E = 0.8396701490447454;  F = 0.9089014466744505 },
    @{ Row = 15; C = "111110101001010011110110110111"; D = 0.9581189667138894; E = 0.8396701490447454;  F = 0.9348716893313952 },
    @{ Row = 16; C = "111110101001010011110111000111"; D = 0.9581189958854615; E = 0.8396701217358291;  F = 0.9157428842300639 },
    @{ Row = 17; C = "111110101001010011110111000111"; D = 0.9581189958854615; E = 0.8396701490447454;  F = 0.9180233681630383 },
    @{ Row = 18; C = "111110101001010011110111000111"; D = 0.9581189958854615; E = 0.8396701490447454;  F = 0.9157428882366567 },
    @{ Row = 19; C = "111110101001010011110111000111"; D = 0.9581189958854615; E = 0.8396701490447454;  F = 0.9061784922304383 },
    @{ Row = 20; C = "111110101001010011110111000111"; D = 0.9581189958854615; E = 0.8396701490447454;  F = 0.8847691312586502 },
    @{ Row = 21; C = "111110101001010011110111000111"; D = 0.9581189958854615; E = 0.8396701490447454;  F = 0.8942969128086302 }
)

foreach ($item in $data) {
    $r = $item.Row

    # Column C holds a chromosome encoded as a string of 0/1 digits. Force
    # the cell to text *before* writing so Excel doesn't reinterpret the
    # all-digit string as a (scientific-notation) number.
    $cellC = $ws.Cells.Item($r, 3)
    $cellC.NumberFormat = "@"
    $cellC.Value = $item.C
    $cellC.NumberFormat = "General"

    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
}
